# Rename sheet name to pdx_model
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("pdx_models")
$ws.Name = "pdx_model"

# Activating the renamed sheet makes it the active/selected tab, matching
# the accompanying view-state changes in the saved workbook (activeTab,
# tabSelected on the sheetViews, and the reset scroll position).
$ws.Activate()
